$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-09-28 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-09-29 Sunday", 2)
$d.Content.Find.Execute("866÷4=216, 2", $true, $false, $false, $false, $false, $true, 1, $false, "573÷7=81, 6", 2)
$d.Content.Find.Execute("816÷7=116, 4", $true, $false, $false, $false, $false, $true, 1, $false, "607÷2=303, 1", 2)
$d.Content.Find.Execute("687÷3=229, 0", $true, $false, $false, $false, $false, $true, 1, $false, "337÷9=37, 4", 2)
$d.Content.Find.Execute("936÷7=133, 5", $true, $false, $false, $false, $false, $true, 1, $false, "405÷2=202, 1", 2)
$d.Content.Find.Execute("754÷2=377, 0", $true, $false, $false, $false, $false, $true, 1, $false, "295÷3=98, 1", 2)
$d.Content.Find.Execute("285÷7=40, 5", $true, $false, $false, $false, $false, $true, 1, $false, "187÷8=23, 3", 2)
$d.Content.Find.Execute("159÷4=39, 3", $true, $false, $false, $false, $false, $true, 1, $false, "854÷8=106, 6", 2)
$d.Content.Find.Execute("375÷6=62, 3", $true, $false, $false, $false, $false, $true, 1, $false, "133÷6=22, 1", 2)
$d.Content.Find.Execute("177÷9=19, 6", $true, $false, $false, $false, $false, $true, 1, $false, "186÷7=26, 4", 2)
$d.Content.Find.Execute("440÷7=62, 6", $true, $false, $false, $false, $false, $true, 1, $false, "145÷5=29, 0", 2)
$d.Content.Find.Execute("934÷6=155, 4", $true, $false, $false, $false, $false, $true, 1, $false, "210÷3=70, 0", 2)
$d.Content.Find.Execute("486÷2=243, 0", $true, $false, $false, $false, $false, $true, 1, $false, "251÷6=41, 5", 2)
$d.Content.Find.Execute("588÷8=73, 4", $true, $false, $false, $false, $false, $true, 1, $false, "706÷9=78, 4", 2)
$d.Content.Find.Execute("870÷7=124, 2", $true, $false, $false, $false, $false, $true, 1, $false, "532÷2=266, 0", 2)
$d.Content.Find.Execute("108÷2=54, 0", $true, $false, $false, $false, $false, $true, 1, $false, "998÷7=142, 4", 2)
$d.Content.Find.Execute("680÷4=170, 0", $true, $false, $false, $false, $false, $true, 1, $false, "838÷7=119, 5", 2)
$d.Content.Find.Execute("717÷4=179, 1", $true, $false, $false, $false, $false, $true, 1, $false, "841÷7=120, 1", 2)
$d.Content.Find.Execute("160÷4=40, 0", $true, $false, $false, $false, $false, $true, 1, $false, "792÷5=158, 2", 2)
$d.Content.Find.Execute("904÷6=150, 4", $true, $false, $false, $false, $false, $true, 1, $false, "230÷6=38, 2", 2)
$d.Content.Find.Execute("826÷9=91, 7", $true, $false, $false, $false, $false, $true, 1, $false, "913÷8=114, 1", 2)
$d.Content.Find.Execute("811÷2=405, 1", $true, $false, $false, $false, $false, $true, 1, $false, "108÷7=15, 3", 2)
$d.Content.Find.Execute("259÷4=64, 3", $true, $false, $false, $false, $false, $true, 1, $false, "211÷4=52, 3", 2)
$d.Content.Find.Execute("991÷3=330, 1", $true, $false, $false, $false, $false, $true, 1, $false, "689÷6=114, 5", 2)
$d.Content.Find.Execute("374÷4=93, 2", $true, $false, $false, $false, $false, $true, 1, $false, "409÷5=81, 4", 2)
$d.Content.Find.Execute("255÷4=63, 3", $true, $false, $false, $false, $false, $true, 1, $false, "955÷7=136, 3", 2)
